$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2235329.8
$ws.Range("J17").Value = 2293988.2
$ws.Range("L17").Value = 6881964.600000001
$ws.Range("N17").Value = -6882300.600000001

$ws.Range("H38").Value = 1403.2
$ws.Range("I38").Value = 404.7647
$ws.Range("J38").Value = 3524.875
$ws.Range("K38").Value = 1214.2941
$ws.Range("L38").Value = 10574.625
$ws.Range("M38").Value = -842.2941000000001
$ws.Range("N38").Value = -11318.625

$ws.Range("H40").Value = 2201.4211
$ws.Range("J40").Value = 2285.2856
$ws.Range("L40").Value = 2285.2856
$ws.Range("N40").Value = -2635.2856

$ws.Range("H64").Value = 6308.478
$ws.Range("I64").Value = 3853.4285
$ws.Range("K64").Value = 3853.4285
$ws.Range("M64").Value = -3605.4285

$ws.Range("H67").Value = 6308.478
$ws.Range("I67").Value = 3853.4285
$ws.Range("K67").Value = 3853.4285
$ws.Range("M67").Value = -2995.4285

$ws.Range("H76").Value = 4064.75
$ws.Range("I76").Value = 4080.182
$ws.Range("K76").Value = 4080.182
$ws.Range("M76").Value = -3765.182

$ws.Range("H79").Value = 4064.75
$ws.Range("I79").Value = 4080.182
$ws.Range("K79").Value = 4080.182
$ws.Range("M79").Value = -2988.182

$ws.Range("H125").Value = 1026.7693
$ws.Range("I125").Value = 854.8570999999999
$ws.Range("K125").Value = 7693.7139
$ws.Range("M125").Value = -5233.7139

$ws.Range("H135").Value = 3707.6428
$ws.Range("I135").Value = 3707.6428
$ws.Range("K135").Value = 33368.7852
$ws.Range("M135").Value = -30833.7852

$ws.Range("H138").Value = 2082.625
$ws.Range("I138").Value = 2276.2942
$ws.Range("J138").Value = 1612.2858
$ws.Range("K138").Value = 6828.882599999999
$ws.Range("L138").Value = 4836.857400000001
$ws.Range("M138").Value = -1688.882599999999
$ws.Range("N138").Value = -15116.8574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3700.8572
$ws.Range("I63").Value = 2180
$ws.Range("J63").Value = 7503
$ws.Range("K63").Value = 2180
$ws.Range("L63").Value = 7503
$ws.Range("M63").Value = -1494
$ws.Range("N63").Value = -8875

$ws.Range("H66").Value = 3700.8572
$ws.Range("I66").Value = 2180
$ws.Range("J66").Value = 7503
$ws.Range("K66").Value = 10900
$ws.Range("L66").Value = 37515
$ws.Range("M66").Value = -7468
$ws.Range("N66").Value = -44379

$ws.Range("H88").Value = 894.53845
$ws.Range("J88").Value = 891.1111
$ws.Range("L88").Value = 891.1111
$ws.Range("N88").Value = -1703.1111

$ws.Range("H91").Value = 894.53845
$ws.Range("J91").Value = 891.1111
$ws.Range("L91").Value = 891.1111
$ws.Range("N91").Value = -3699.1111

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 2636.44
$ws.Range("I122").Value = 2619.2632
$ws.Range("K122").Value = 7857.7896
$ws.Range("M122").Value = -5407.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5408.5293
$ws.Range("I105").Value = 5093.7
$ws.Range("K105").Value = 5093.7
$ws.Range("M105").Value = -3346.7

$ws.Range("H107").Value = 2133.5417
$ws.Range("I107").Value = 1622.9
$ws.Range("K107").Value = 1622.9
$ws.Range("M107").Value = 297.0999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5676.65
$ws.Range("I58").Value = 3651.6
$ws.Range("K58").Value = 3651.6
$ws.Range("M58").Value = -3448.6

$ws.Range("H62").Value = 89729.914
$ws.Range("I62").Value = 169548.17
$ws.Range("J62").Value = 9911.666999999999
$ws.Range("K62").Value = 169548.17
$ws.Range("L62").Value = 9911.666999999999
$ws.Range("M62").Value = -168924.17
$ws.Range("N62").Value = -11159.667

$ws.Range("H65").Value = 89729.914
$ws.Range("I65").Value = 169548.17
$ws.Range("J65").Value = 9911.666999999999
$ws.Range("K65").Value = 847740.8500000001
$ws.Range("L65").Value = 49558.335
$ws.Range("M65").Value = -844620.8500000001
$ws.Range("N65").Value = -55798.335

$ws.Range("H132").Value = 33346.28
$ws.Range("I132").Value = 24887.285
$ws.Range("K132").Value = 74661.855
$ws.Range("M132").Value = -72131.855

$ws.Range("H136").Value = 5676.65
$ws.Range("I136").Value = 3651.6
$ws.Range("K136").Value = 10954.8
$ws.Range("M136").Value = -8404.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 2331.6667
$ws.Range("J29").Value = 2331.6667
$ws.Range("L29").Value = 6995.000100000001
$ws.Range("N29").Value = -7549.000100000001

$ws.Range("H34").Value = 2695.742
$ws.Range("I34").Value = 171.375
$ws.Range("J34").Value = 5388.4
$ws.Range("K34").Value = 514.125
$ws.Range("L34").Value = 16165.2
$ws.Range("M34").Value = -430.125
$ws.Range("N34").Value = -16333.2

$ws.Range("H131").Value = 2995.7856
$ws.Range("J131").Value = 4269
$ws.Range("L131").Value = 12807
$ws.Range("N131").Value = -22887

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7797.85
$ws.Range("I70").Value = 7373.625
$ws.Range("J70").Value = 8080.6665
$ws.Range("K70").Value = 7373.625
$ws.Range("L70").Value = 8080.6665
$ws.Range("M70").Value = -7103.625
$ws.Range("N70").Value = -8620.666499999999

$ws.Range("H73").Value = 7797.85
$ws.Range("I73").Value = 7373.625
$ws.Range("J73").Value = 8080.6665
$ws.Range("K73").Value = 7373.625
$ws.Range("L73").Value = 8080.6665
$ws.Range("M73").Value = -6437.625
$ws.Range("N73").Value = -9952.666499999999

$ws.Range("I80").Value = 2098.3333
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2098.3333
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1100.3333
$ws.Range("N80").ClearContents()

$ws.Range("I83").Value = 2098.3333
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 10491.6665
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5499.666499999999
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 23366.965
$ws.Range("I132").Value = 10881.521
$ws.Range("K132").Value = 32644.563
$ws.Range("M132").Value = -30114.563

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2040.1904
$ws.Range("I46").Value = 1228.5714
$ws.Range("J46").Value = 2446
$ws.Range("K46").Value = 1228.5714
$ws.Range("L46").Value = 2446
$ws.Range("M46").Value = -1040.5714
$ws.Range("N46").Value = -2822

$ws.Range("H136").Value = 5735.0938
$ws.Range("I136").Value = 4946.1577
$ws.Range("K136").Value = 14838.4731
$ws.Range("M136").Value = -12288.4731

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 24150.5
$ws.Range("J103").Value = 24150.5
$ws.Range("L103").Value = 24150.5
$ws.Range("N103").Value = -26494.5

$ws.Range("H113").Value = 624.4138
$ws.Range("I113").Value = 503.89285
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 1511.67855
$ws.Range("L113").Value = 11997
$ws.Range("M113").Value = 658.3214499999999
$ws.Range("N113").Value = -16337

$ws.Range("H136").Value = 6898855
$ws.Range("I136").Value = 11766959
$ws.Range("K136").Value = 35300877
$ws.Range("M136").Value = -35298327
